# Apply changes described in commit: "added institut ID and case ID to DQ Reports"

$wb = $excel.ActiveWorkbook

# --- Sheet: DQ_Report ---
# Insert a new column before the existing "ICD_primaerkode" column (column B)
# to hold the new "Aufnahmenummer" (case/admission ID) values.
$wsDQ = $wb.Worksheets.Item("DQ_Report")
$wsDQ.Range("B1").EntireColumn.Insert()

# Header for the newly inserted column
$wsDQ.Cells.Item(1, 2).Value = "Aufnahmenummer"

# Fill in the new "Aufnahmenummer" values for each data row
$wsDQ.Cells.Item(2, 2).Value = "F_101645"
$wsDQ.Cells.Item(3, 2).Value = "F_101646"
$wsDQ.Cells.Item(4, 2).Value = "F_101648"
$wsDQ.Cells.Item(5, 2).Value = "F_101649"
$wsDQ.Cells.Item(6, 2).Value = "F_101650"
$wsDQ.Cells.Item(7, 2).Value = "F_101651"
$wsDQ.Cells.Item(8, 2).Value = "F_101651"
$wsDQ.Cells.Item(9, 2).Value = "F_101653"
$wsDQ.Cells.Item(10, 2).Value = "F_101654"
$wsDQ.Cells.Item(11, 2).Value = "F_101655"
$wsDQ.Cells.Item(12, 2).Value = "F_101656"
$wsDQ.Cells.Item(13, 2).Value = "F_101757"
$wsDQ.Cells.Item(14, 2).Value = "F_101658"
$wsDQ.Cells.Item(15, 2).Value = "F_101660"

# --- Sheet: Statistik ---
# Rename the "basicItem" column header to "inst_id" and populate the
# institute ID / updated statistics for the first (and only) data row.
$wsStat = $wb.Worksheets.Item("Statistik")
$wsStat.Cells.Item(1, 1).Value = "inst_id"
$wsStat.Cells.Item(2, 1).Value = "260123430-Dali"
$wsStat.Cells.Item(2, 2).Value = 3.13
$wsStat.Cells.Item(2, 3).Value = 96.87
